$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Column B (Coin name) updates (rows 20/21 swap) ---
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('B21').Value = 'Uniswap'

# --- Column C (Link) updates (rows 20/21 swap) ---
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'

# --- Column D (Price) updates - forced to remain text, matching original inlineStr type ---
Set-TextValue $ws 'D2' '26.227.65'
Set-TextValue $ws 'D3' '1.645.53'
Set-TextValue $ws 'D5' '217.26'
Set-TextValue $ws 'D9' '0.0638'
Set-TextValue $ws 'D10' '20.01'
Set-TextValue $ws 'D13' '1.873.30'
Set-TextValue $ws 'D14' '1.663.84'
Set-TextValue $ws 'D16' '0.0₃0765'
Set-TextValue $ws 'D18' '26.216.74'
Set-TextValue $ws 'D20' '195.54'
Set-TextValue $ws 'D21' '4.44'
Set-TextValue $ws 'D23' '6.34'
Set-TextValue $ws 'D24' '143.48'
Set-TextValue $ws 'D26' '1.77'
Set-TextValue $ws 'D28' '6.95'
Set-TextValue $ws 'D29' '15.63'
Set-TextValue $ws 'D32' '3.35'
Set-TextValue $ws 'D34' '1.60'
Set-TextValue $ws 'D38' '1.133.10'
Set-TextValue $ws 'D43' '100.12'
Set-TextValue $ws 'D45' '1.782.74'
Set-TextValue $ws 'D46' '56.31'
Set-TextValue $ws 'D49' '7.72'
Set-TextValue $ws 'D51' '0.0974'

# --- Column E (Volume 1h) updates - forced to remain text, matching original inlineStr type ---
Set-TextValue $ws 'E2' '  +1.64%  '
Set-TextValue $ws 'E3' '  +0.55%  '
Set-TextValue $ws 'E4' '  -0.15%  '
Set-TextValue $ws 'E5' '  +0.82%  '
Set-TextValue $ws 'E7' '  -0.15%  '
Set-TextValue $ws 'E8' '  +0.10%  '
Set-TextValue $ws 'E9' '  +0.11%  '
Set-TextValue $ws 'E10' '  +1.35%  '
Set-TextValue $ws 'E11' '  +0.03%  '
Set-TextValue $ws 'E12' '  +0.73%  '
Set-TextValue $ws 'E13' '  +0.59%  '
Set-TextValue $ws 'E14' '  +1.75%  '
Set-TextValue $ws 'E15' '  -2.24%  '
Set-TextValue $ws 'E16' '  -0.31%  '
Set-TextValue $ws 'E17' '  +0.53%  '
Set-TextValue $ws 'E18' '  +1.52%  '
Set-TextValue $ws 'E19' '  -0.13%  '
Set-TextValue $ws 'E20' '  +1.56%  '
Set-TextValue $ws 'E21' '  -0.63%  '
Set-TextValue $ws 'E22' '  +0.74%  '
Set-TextValue $ws 'E23' '  -0.37%  '
Set-TextValue $ws 'E24' '  +0.75%  '
Set-TextValue $ws 'E25' '  -0.15%  '
Set-TextValue $ws 'E26' '  -2.60%  '
Set-TextValue $ws 'E27' '  +1.84%  '
Set-TextValue $ws 'E28' '  -0.12%  '
Set-TextValue $ws 'E29' '  +0.52%  '
Set-TextValue $ws 'E30' '  +1.36%  '
Set-TextValue $ws 'E31' '  +1.98%  '
Set-TextValue $ws 'E32' '  +0.05%  '
Set-TextValue $ws 'E34' '  +1.58%  '
Set-TextValue $ws 'E35' '  +0.79%  '
Set-TextValue $ws 'E37' '  +1.87%  '
Set-TextValue $ws 'E38' '  +0.15%  '
Set-TextValue $ws 'E39' '  -1.68%  '
Set-TextValue $ws 'E40' '  +0.63%  '
Set-TextValue $ws 'E41' '  -0.03%  '
Set-TextValue $ws 'E42' '  +1.84%  '
Set-TextValue $ws 'E43' '  -0.46%  '
Set-TextValue $ws 'E44' '  -1.11%  '
Set-TextValue $ws 'E45' '  +0.61%  '
Set-TextValue $ws 'E46' '  +1.78%  '
Set-TextValue $ws 'E47' '  +4.60%  '
Set-TextValue $ws 'E48' '  +2.62%  '
Set-TextValue $ws 'E49' '  +3.29%  '
Set-TextValue $ws 'E50' '  +0.05%  '
Set-TextValue $ws 'E51' '  +1.85%  '

